$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G3").Value = "servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg , hend_mahmoud@med.asu.edu.eg"
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G14").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G27").Value = "randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G29").Value = "aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G33").Value = "servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G34").Value = "servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg , hend_mahmoud@med.asu.edu.eg"
$ws.Range("G40").Value = "Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G45").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G46").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G58").Value = "randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G60").Value = "aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G64").Value = "mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G65").Value = "Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G71").Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G76").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G77").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G91").Value = "aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G95").Value = "servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G96").Value = "servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, nourhan.mahmoud@med.asu.edu.eg "
$ws.Range("G102").Value = "Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G107").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G108").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G120").Value = "marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg"
$ws.Range("G122").Value = "Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G126").Value = "servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G127").Value = "servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, nourhan.mahmoud@med.asu.edu.eg "
$ws.Range("G133").Value = "Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G138").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G139").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G140").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G151").Value = "user@user.com, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg`n"
$ws.Range("G153").Value = "Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G157").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G158").Value = "Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G164").Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G169").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G170").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G182").Value = "marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg"
$ws.Range("G184").Value = "maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
